# Update cryptocurrency price/volume data (and reorder three coin rows)
# to reflect the refreshed GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.142.65"
$ws.Range("E2").Value = "  +12.47%  "

$ws.Range("D3").Value = "1.619.61"
$ws.Range("E3").Value = "  +10.09%  "

$ws.Range("D4").Value = "'0.9886"
$ws.Range("E4").Value = "  -1.94%  "

$ws.Range("D5").Value = "'301.40"
$ws.Range("E5").Value = "  +8.78%  "

$ws.Range("D6").Value = "'0.9814"
$ws.Range("E6").Value = "  +2.51%  "

$ws.Range("D7").Value = "'0.3662"
$ws.Range("E7").Value = "  +2.81%  "

$ws.Range("D8").Value = "'0.3414"
$ws.Range("E8").Value = "  +11.25%  "

$ws.Range("D9").Value = "'42.13"
$ws.Range("E9").Value = "  +6.92%  "

$ws.Range("D10").Value = "'1.140"
$ws.Range("E10").Value = "  +4.79%  "

$ws.Range("D11").Value = "'0.07053"
$ws.Range("E11").Value = "  +6.31%  "

$ws.Range("D12").Value = "'0.9850"
$ws.Range("E12").Value = "  -1.72%  "

$ws.Range("D13").Value = "'19.97"
$ws.Range("E13").Value = "  +10.36%  "

$ws.Range("D14").Value = "'5.877"
$ws.Range("E14").Value = "  +7.54%  "

$ws.Range("D15").Value = "'6.595"
$ws.Range("E15").Value = "  +6.85%  "

$ws.Range("D16").Value = "'0.00001078"
$ws.Range("E16").Value = "  +5.48%  "

$ws.Range("D17").Value = "1.619.05"
$ws.Range("E17").Value = "  +10.27%  "

$ws.Range("D18").Value = "'0.9811"
$ws.Range("E18").Value = "  +2.45%  "

$ws.Range("D19").Value = "'0.06646"
$ws.Range("E19").Value = "  +11.32%  "

$ws.Range("D20").Value = "'78.29"
$ws.Range("E20").Value = "  +13.45%  "

$ws.Range("D21").Value = "'16.08"
$ws.Range("E21").Value = "  +10.62%  "

$ws.Range("D22").Value = "'5.981"
$ws.Range("E22").Value = "  +9.12%  "

$ws.Range("D23").Value = "'11.67"
$ws.Range("E23").Value = "  +3.42%  "

$ws.Range("D24").Value = "23.088.33"
$ws.Range("E24").Value = "  +12.28%  "

$ws.Range("D25").Value = "'2.376"
$ws.Range("E25").Value = "  +4.89%  "

$ws.Range("D26").Value = "'2.595"
$ws.Range("E26").Value = "  +24.12%  "

$ws.Range("D27").Value = "'149.89"
$ws.Range("E27").Value = "  +3.26%  "

$ws.Range("D28").Value = "'19.32"
$ws.Range("E28").Value = "  +12.86%  "

$ws.Range("D29").Value = "1.794.00"
$ws.Range("E29").Value = "  +10.18%  "

$ws.Range("D30").Value = "'124.81"
$ws.Range("E30").Value = "  +9.52%  "

$ws.Range("D31").Value = "'4.074"
$ws.Range("E31").Value = "  +5.62%  "

$ws.Range("D32").Value = "'6.002"
$ws.Range("E32").Value = "  +21.64%  "

$ws.Range("D33").Value = "'0.9717"
$ws.Range("E33").Value = "  +21.78%  "

$ws.Range("D34").Value = "'1.662"
$ws.Range("E34").Value = "  +15.39%  "

$ws.Range("D35").Value = "'0.08216"
$ws.Range("E35").Value = "  +3.52%  "

$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").Value = "'8.701"
$ws.Range("E36").Value = "  +19.62%  "

$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").Value = "'11.77"
$ws.Range("E37").Value = "  +13.85%  "

$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").Value = "'5.157"
$ws.Range("E38").Value = "  +9.40%  "

$ws.Range("D39").Value = "'1.252"
$ws.Range("E39").Value = "  +0.87%  "

$ws.Range("D40").Value = "'0.06121"
$ws.Range("E40").Value = "  +6.62%  "

$ws.Range("D41").Value = "'0.02212"
$ws.Range("E41").Value = "  +9.00%  "

$ws.Range("D42").Value = "'0.2014"
$ws.Range("E42").Value = "  +8.21%  "

$ws.Range("D43").Value = "'0.9811"
$ws.Range("E43").Value = "  +2.41%  "

$ws.Range("D44").Value = "'0.5871"
$ws.Range("E44").Value = "  +11.76%  "

$ws.Range("D45").Value = "'3.774"
$ws.Range("E45").Value = "  +7.45%  "

$ws.Range("D46").Value = "'13.00"
$ws.Range("E46").Value = "  +7.27%  "

$ws.Range("D47").Value = "'0.5730"
$ws.Range("E47").Value = "  +10.58%  "

$ws.Range("D48").Value = "'125.97"
$ws.Range("E48").Value = "  +6.21%  "

$ws.Range("E49").Value = "  +9.22%  "

$ws.Range("D50").Value = "'0.06932"
$ws.Range("E50").Value = "  +7.74%  "

$ws.Range("D51").Value = "'73.55"
$ws.Range("E51").Value = "  +9.92%  "
